$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header labels in row 1 need to be capitalized so the country /
# country_de columns are correctly recognized (picked up) downstream.
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Country_de"

$ws.Range("B2").Select()
